$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Cypher query text for the "ParticipantsTab" row (cell B2) - replaces
# the previous query with an updated one that walks optional matches and
# sorts the collected sample ids.
$newQuery = "MATCH (p:participant)-->(s:study)" + "`n" + `
"OPTIONAL MATCH (samp:sample)-->(p)" + "`n" + `
"OPTIONAL MATCH (p)<--(diag:diagnosis)" + "`n" + `
"OPTIONAL MATCH (samp)<--(f:file)" + "`n" + `
"OPTIONAL MATCH (f)<--(g:genomic_info)" + "`n" + `
"WITH s, p, samp, f, g, diag" + "`n" + `
"WHERE g.platform in ['NovaSeqS4']" + "`n" + `
"with p" + "`n" + `
"OPTIONAL MATCH (p)-->(s:study)" + "`n" + `
"OPTIONAL MATCH (samp:sample)-->(p)" + "`n" + `
"WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp" + "`n" + `
"RETURN" + "`n" + `
"coalesce(p.participant_id,'') as ``Participant ID``," + "`n" + `
"coalesce(s.study_name, '') as ``Study Name``," + "`n" + `
"coalesce(s.phs_accession,'') as ``Accession``," + "`n" + `
"coalesce(p.gender,'') as ``Gender``," + "`n" + `
"coalesce(apoc.text.join(samp, ','), '') as ``Samples``" + "`n" + `
"ORDER BY p.participant_id LIMIT 100"

$ws.Range("B2").Value = $newQuery

# The replacement query text wraps to more lines than the original, so the
# row grows taller to keep the whole query visible.
$ws.Rows(2).RowHeight = 279

# Move the active selection to B4.
$ws.Range("B4").Select()
